$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 4)
# Order of assignment matches the order new strings were appended to the
# shared string table: Susquehanna, job-location, https://sig.com/, careers URL
$ws.Range("A4").Value = "Susquehanna"
$ws.Range("E4").Value = "job-location"
$ws.Range("C4").Value = "https://sig.com/"
$ws.Range("B4").Value = "https://careers.sig.com/c/quantitative-trading-strategy-jobs"
$ws.Range("D4").Value = "job-title"

# Update the selection to match the author's final cursor position
$ws.Range("B8").Select()
